$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 72 (between "Charge Sequence" and "Charge Description") for "Charge Severity"
$ws.Rows("72:72").Insert(-4121, 0)
$ws.Range("B72").Value = "Charge Severity"
$ws.Range("C72").Value = "Charge Severity Text"
$ws.Range("E72").Value = "/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Charge[@structures:id=/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Arrest/j:ArrestCharge/@structures:ref]/j:ChargeSeverityText"

# Insert a new row at 76 (between "Charge Category/Classification" and "Highest Charge Category") for "Charge Jurisdiction Court"
$ws.Rows("76:76").Insert(-4121, 0)
$ws.Range("A76").Value = $ws.Range("A77").Value2
$ws.Range("A77").Value = ""
$ws.Range("B76").Value = "Charge Jurisdiction Court"
$ws.Range("E76").Value = "/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Charge[@structures:id=/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Arrest/j:ArrestCharge/@structures:ref]/cscr-ext:ChargeJurisdictionCourt/j:CourtName"

# Update view pane/selection
$ws.Application.ActiveWindow.ScrollRow = 66
$ws.Range("E76").Select()
